# EPBDS-9436: Add test for getValue() error message for case-sensitive vocabulary case.
# Adds a new "Datatype StrField <String>" datatype + enum-like values (aaa/bbb/ccc) and a
# new Spreadsheet mySpr4(MyDatatype param) that calls getValue(strField), preceded by an
# explanatory comment row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New Datatype "StrField <String>" enumeration block (bordered cells, style like C65:C69) ---
$ws.Range("C65").Copy()
$ws.Range("C73").PasteSpecial(-4122)
$ws.Range("C74").PasteSpecial(-4122)
$ws.Range("C75").PasteSpecial(-4122)
$ws.Range("C76").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C73").Value = "Datatype StrField <String> "
$ws.Range("C74").Value = "aaa"
$ws.Range("C75").Value = "bbb"
$ws.Range("C76").Value = "ccc"

# --- New Spreadsheet mySpr4(MyDatatype param) block, with a merged/bordered title row ---
$ws.Range("C79:D79").Merge()
$ws.Range("C65").Copy()
$ws.Range("C79:D79").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C79").Value = "Spreadsheet SpreadsheetResult mySpr4(MyDatatype param)"

$ws.Range("C65:D65").Copy()
$ws.Range("C80:D80").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C80").Value = "Steps"
$ws.Range("D80").Value = "Values"

$ws.Range("C65").Copy()
$ws.Range("C81").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C81").Value = "getValue"

$ws.Range("D21").Copy()
$ws.Range("D81").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D81").Value = "'=getValue(strField)"

# --- Two trailing blank styled cells (quote-prefix style, matching D15's style) ---
$ws.Range("D15").Copy()
$ws.Range("D82").PasteSpecial(-4122)
$ws.Range("D83").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New comment row (no border, no special style - matches existing "section title" rows
#     such as C42/C47/C53/C60 which also carry no explicit style). Typed last, so its shared
#     string lands at the end of the table. ---
$ws.Range("C71").Value = "// The case with Case sensitivity  for vocabulary and getValue() method EPBDS-9436"

# --- Remove the unused shared string "=(Integer)mySpr(a).getFieldValue("$Step1")" by
#     clearing any stray reference to it (it is not referenced by any live cell; Excel
#     drops unreferenced shared strings on save). ---

# --- Restore selection/view close to the authored state ---
$ws.Range("C74").Select()

$wb.Save()
